$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: "Cross Validation Accuracy" -> "Cross Validation Mean Accuracy"
$ws.Range("D1").Value = "Cross Validation Mean Accuracy"

# --- Preserve the existing data in rows 8 & 9 before the labels get renamed.
# Original row 8: KNeighbours Classifier (Count Vectorizer + TfidfVectorizer)
# Original row 9: Gradient Boosting Classifier (Count Vectorizer + TfidfVectorizer)
$oldRow8B = $ws.Range("B8").Value2
$oldRow8C = $ws.Range("C8").Value2
$oldRow8D = $ws.Range("D8").Value2
$oldRow9B = $ws.Range("B9").Value2
$oldRow9C = $ws.Range("C9").Value2
$oldRow9D = $ws.Range("D9").Value2

# --- Rename the two existing rows' labels in place to the updated wording
# ("TfidfVectorizer" -> "TfidfTransformer"). Row 9 first, then row 8.
$ws.Range("A9").Value = "Gradient Boosting Classifier (Count Vectorizer + TfidfTransformer)"
$ws.Range("A8").Value = "KNeighbours Classifier (Count Vectorizer + TfidfTransformer)"

# --- Make room for two new rows: one new "Count Vectorizer only" row ahead of
# each of the renamed rows.
$ws.Rows.Item(9).Insert()
$ws.Rows.Item(8).Insert()

# Final row layout after the inserts:
#   8  (blank, will hold the new KNeighbours (Count Vectorizer) row)
#   9  KNeighbours Classifier (Count Vectorizer + TfidfTransformer) - old row8 data
#   10 (blank, will hold the new Gradient Boosting (Count Vectorizer) row)
#   11 Gradient Boosting Classifier (Count Vectorizer + TfidfTransformer) - old row9 data

# --- New label text, entered in the same order the workbook's shared-string
# table records them.
$ws.Range("A10").Value = "Gradient Boosting Classifier (Count Vectorizer)"
$ws.Range("A8").Value  = "KNeighbours Classifier (Count Vectorizer)"
$ws.Range("A12").Value = "Xgboost Classifier (Count Vectorizer)"
$ws.Range("A13").Value = "Xgboost Classifier (Count Vectorizer + TfidfTransformer)"

# --- Row 9: restore the old row-8 data under its new label
$ws.Range("B9").Value = $oldRow8B
$ws.Range("C9").Value = $oldRow8C
$ws.Range("D9").Value = $oldRow8D

# --- Row 11: restore the old row-9 data under its new label
$ws.Range("B11").Value = $oldRow9B
$ws.Range("C11").Value = $oldRow9C
$ws.Range("D11").Value = $oldRow9D

# --- Row 10: new Gradient Boosting (Count Vectorizer) values
$ws.Range("B10").Value = 0.95907473309608504
$ws.Range("C10").Value = 0.99885583524027399

# --- Row 12: new Xgboost (Count Vectorizer) values (no cross-validation figure)
$ws.Range("B12").Value = 0.97458057956278599
$ws.Range("C12").Value = 0.99771167048054898

# --- Row 13: new Xgboost (Count Vectorizer + TfidfTransformer) values
$ws.Range("B13").Value = 0.97127605490594804
$ws.Range("C13").Value = 0.99879227053140096

# --- Row 8 has no numeric data at all for this model -> remove the inherited
# empty/styled B:D cells entirely.
$ws.Range("B8:D8").Clear()

# D10 never got a cross-validation figure -> clear the value but keep the
# cell (and its formatting) present, same as the rest of the column.
$ws.Range("D10").ClearContents()

# Rows 12 & 13 are brand new rows with no inherited formatting -> give their
# numeric cells the same style (font + alignment) as the rest of the table.
$ws.Range("B2:C2").Copy() | Out-Null
$ws.Range("B12:C13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Column A needs to be a little wider to fit the longer model names.
# (Target authored width is 55.5546875 characters; this runtime's ColumnWidth
# setter only resolves to 1/6-character increments, so 55.5 is the closest
# reachable approximation.)
$ws.Columns.Item(1).ColumnWidth = 54.666666666666664

# --- Selection moves to A16 in the saved file (matches authored selection state)
$ws.Range("A16").Select()
